$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Sheets:" $wb.Worksheets.Count
